# Applies the authoritative edit: cyclic rotation of observation rows 3-6.
# Row 3 <- old row 4, Row 4 <- old row 5, Row 5 <- old row 6, Row 6 <- old row 3
# (each full record moves up one row; the Knärot record from row 3 wraps to row 6).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 ---
$ws.Range("A3").Value = 130937843
$ws.Range("B3").Value = 57884
$ws.Range("D3").Value = "NT"
$ws.Range("E3").Value = 100109
$ws.Range("F3").Value = "Tretåig hackspett"
$ws.Range("G3").Value = "Picoides tridactylus"
$ws.Range("H3").Value = "(Linnaeus, 1758)"
$ws.Range("I3").Value = ""
$ws.Range("J3").ClearContents()
$ws.Range("K3").Value = ""
$ws.Range("M3").Value = "färska spår"
$ws.Range("P3").Value = "Storflon, Jmt"
$ws.Range("Q3").Value = 489760
$ws.Range("R3").Value = 7004232
$ws.Range("S3").Value = 10
$ws.Range("T3").Value = "Jämtland"
$ws.Range("U3").Value = "Östersund"
$ws.Range("V3").Value = "Jämtland"
$ws.Range("W3").Value = "Brunflo"
$ws.Range("Y3").NumberFormat = "@"
$ws.Range("Y3").Value = "2026-01-28"
$ws.Range("Y3").ClearFormats()
$ws.Range("AA3").NumberFormat = "@"
$ws.Range("AA3").Value = "2026-01-28"
$ws.Range("AA3").ClearFormats()
$ws.Range("AC3").Value = "Ringhack, färska och äldre, i riklig mängd längs flera meter högt upp på en granstam med spår av rikligt sav/kådaflöde."
$ws.Range("AD3").Value = $False
$ws.Range("AE3").Value = $False
$ws.Range("AF3").ClearContents()
$ws.Range("AG3").Value = $False
$ws.Range("AH3").Value = "Granskog"
$ws.Range("AJ3").Value = "gran"
$ws.Range("AK3").Value = "Picea abies"
$ws.Range("AM3").Value = "Trädstam på levande träd"
$ws.Range("AO3").Value = "Stem on living tree # Picea abies"
$ws.Range("AW3").Value = "Kristian Zackrisson"
$ws.Range("AX3").Value = "Kristian Zackrisson"

# --- Row 4 ---
$ws.Range("A4").Value = 130937852
$ws.Range("B4").Value = 57884
$ws.Range("D4").Value = "NT"
$ws.Range("E4").Value = 100109
$ws.Range("F4").Value = "Tretåig hackspett"
$ws.Range("G4").Value = "Picoides tridactylus"
$ws.Range("H4").Value = "(Linnaeus, 1758)"
$ws.Range("I4").Value = ""
$ws.Range("J4").ClearContents()
$ws.Range("K4").Value = ""
$ws.Range("M4").Value = "äldre spår"
$ws.Range("P4").Value = "Storflon, Jmt"
$ws.Range("Q4").Value = 489520
$ws.Range("R4").Value = 7004161
$ws.Range("S4").Value = 10
$ws.Range("T4").Value = "Jämtland"
$ws.Range("U4").Value = "Östersund"
$ws.Range("V4").Value = "Jämtland"
$ws.Range("W4").Value = "Brunflo"
$ws.Range("Y4").NumberFormat = "@"
$ws.Range("Y4").Value = "2026-01-28"
$ws.Range("Y4").ClearFormats()
$ws.Range("AA4").NumberFormat = "@"
$ws.Range("AA4").Value = "2026-01-28"
$ws.Range("AA4").ClearFormats()
$ws.Range("AC4").Value = "Ringhack, äldre, ytliga enstaka längs flera meter på en granstam vid kanten mot yngre skog."
$ws.Range("AD4").Value = $False
$ws.Range("AE4").Value = $False
$ws.Range("AF4").ClearContents()
$ws.Range("AG4").Value = $False
$ws.Range("AH4").Value = "Granskog"
$ws.Range("AJ4").Value = "gran"
$ws.Range("AK4").Value = "Picea abies"
$ws.Range("AM4").Value = "Trädstam på levande träd"
$ws.Range("AO4").Value = "Stem on living tree # Picea abies"
$ws.Range("AW4").Value = "Kristian Zackrisson"
$ws.Range("AX4").Value = "Kristian Zackrisson"

# --- Row 5 ---
$ws.Range("A5").Value = 130937854
$ws.Range("B5").Value = 57881
$ws.Range("D5").Value = "NT"
$ws.Range("E5").Value = 100049
$ws.Range("F5").Value = "Spillkråka"
$ws.Range("G5").Value = "Dryocopus martius"
$ws.Range("H5").Value = "(Linnaeus, 1758)"
$ws.Range("I5").Value = ""
$ws.Range("J5").ClearContents()
$ws.Range("K5").Value = ""
$ws.Range("M5").Value = "färska spår"
$ws.Range("P5").Value = "Storflon, Jmt"
$ws.Range("Q5").Value = 489668
$ws.Range("R5").Value = 7004128
$ws.Range("S5").Value = 10
$ws.Range("T5").Value = "Jämtland"
$ws.Range("U5").Value = "Östersund"
$ws.Range("V5").Value = "Jämtland"
$ws.Range("W5").Value = "Brunflo"
$ws.Range("Y5").NumberFormat = "@"
$ws.Range("Y5").Value = "2026-01-28"
$ws.Range("Y5").ClearFormats()
$ws.Range("AA5").NumberFormat = "@"
$ws.Range("AA5").Value = "2026-01-28"
$ws.Range("AA5").ClearFormats()
$ws.Range("AC5").Value = "Rejäla hackspår, färska och äldre, I två levande granar och i ytlig grov rotdel."
$ws.Range("AD5").Value = $False
$ws.Range("AE5").Value = $False
$ws.Range("AF5").ClearContents()
$ws.Range("AG5").Value = $False
$ws.Range("AH5").Value = "Granskog"
$ws.Range("AJ5").Value = "gran"
$ws.Range("AK5").Value = "Picea abies"
$ws.Range("AM5").Value = "Trädstam på levande träd"
$ws.Range("AO5").Value = "Stem on living tree # Picea abies"
$ws.Range("AW5").Value = "Kristian Zackrisson"
$ws.Range("AX5").Value = "Kristian Zackrisson"

# --- Row 6 ---
$ws.Range("A6").Value = 130937863
$ws.Range("B6").Value = 99015
$ws.Range("D6").Value = "VU"
$ws.Range("E6").Value = 220787
$ws.Range("F6").Value = "Knärot"
$ws.Range("G6").Value = "Goodyera repens"
$ws.Range("H6").Value = "(L.) R. Br."
$ws.Range("I6").NumberFormat = "@"
$ws.Range("I6").Value = "8"
$ws.Range("I6").ClearFormats()
$ws.Range("J6").Value = "plantor/tuvor"
$ws.Range("K6").Value = "fullt utvecklade blad"
$ws.Range("M6").ClearContents()
$ws.Range("P6").Value = "Storflon, Jmt"
$ws.Range("Q6").Value = 489799
$ws.Range("R6").Value = 7004245
$ws.Range("S6").Value = 10
$ws.Range("T6").Value = "Jämtland"
$ws.Range("U6").Value = "Östersund"
$ws.Range("V6").Value = "Jämtland"
$ws.Range("W6").Value = "Brunflo"
$ws.Range("Y6").NumberFormat = "@"
$ws.Range("Y6").Value = "2026-01-28"
$ws.Range("Y6").ClearFormats()
$ws.Range("AA6").NumberFormat = "@"
$ws.Range("AA6").Value = "2026-01-28"
$ws.Range("AA6").ClearFormats()
$ws.Range("AC6").Value = "Minst 8 plantor inom ca 1 m2 yta. Grävdes varsamt fram under snötäcket. Det finns sannolikt betydligt mer knärot på fyndplatsen och i skogsbeståndet där fyndplatsen ligger."
$ws.Range("AD6").Value = $False
$ws.Range("AE6").Value = $False
$ws.Range("AF6").Value = ""
$ws.Range("AG6").Value = $False
$ws.Range("AH6").Value = "Barrskog"
$ws.Range("AJ6").ClearContents()
$ws.Range("AK6").ClearContents()
$ws.Range("AM6").ClearContents()
$ws.Range("AO6").ClearContents()
$ws.Range("AW6").Value = "Kristian Zackrisson"
$ws.Range("AX6").Value = "Kristian Zackrisson"
